# Error Calculations and Plots
# Apply the edits described by the commit: two sample rows ("RM 232" and
# "SC 92") are removed from the missing-data table, and a number of
# individual cells are toggled between a value and "missing" (blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that were dropped from the dataset -----------------
# Delete from the bottom up so row numbers above the deletion point stay
# valid for the next delete call.
$ws.Rows(28).Delete()   # "SC 92"
$ws.Rows(26).Delete()   # "RM 232"

# --- Individual cell edits (row numbers below refer to the sheet AFTER the
# two rows above have been removed) ------------------------------------------
$ws.Range("C3").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("F6").Value = 16.43
$ws.Range("D8").Value = -13.9
$ws.Range("D10").Value = -14.7
$ws.Range("F11").Value = 17.65
$ws.Range("D12").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("F13").Value = 17.1
$ws.Range("D15").Value = -15.2
$ws.Range("F17").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("D19").ClearContents()
$ws.Range("F19").ClearContents()
$ws.Range("D25").Value = -15.5
$ws.Range("F25").Value = 16.6
$ws.Range("B26").Value = -20.2
$ws.Range("B27").ClearContents()
$ws.Range("D29").ClearContents()
$ws.Range("F31").ClearContents()
$ws.Range("F32").ClearContents()
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
